# add basic xlsx support
# - rename the (only) worksheet from "Лист1" to "1"
# - add a new column C filled with "A" for every existing data row (1-30)
# - scroll the sheet view so that row 17 is the first visible row (best effort)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet "Лист1" -> "1"
$ws.Name = "1"

# Determine the current extent of the data (last used row) so this keeps
# working even if the sheet layout changes.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

# Fill the new column C ("A" for every row) using a single range assignment
# so all the new cells share one shared-string entry.
$ws.Range("C1:C" + $lastRow).Value = "A"

# Best-effort: move the view so row 17 becomes the top-left visible cell
# (mirrors the topLeftCell="A17" change on the sheetView).
$ws.Range("A17").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1

# Restore the original selection (A1), matching the unchanged <selection/>.
$ws.Range("A1").Select()
